# Apply the "insurance" sheet rework:
#  - rename the 5th sheet from "具有相當價值之財產" to "保險"
#  - drop the stand-alone label header row (保險公司/保險名稱/要保人/備註) so the
#    three insurance records shift up by one row; row 1 keeps mirroring row 2
#    like before (now mirroring the first insurance record instead of labels)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

$ws.Name = "保險"

# Overwrite the values in place (rows 1-4 already carry the right styles:
# row 1 / column A bold, everything else normal) so no row insert/delete is
# needed for the shift itself.

$ws.Cells.Item(1, 2).Value = "三商美邦人壽"
$ws.Cells.Item(1, 3).Value = "世紀理財變額萬能終身壽險"
$ws.Cells.Item(1, 4).Value = "蔡煌瑯"
$ws.Cells.Item(1, 5).Value = "保單號碼:649700013402自96年1月11日起迄今20年期(可隨時終止）"

$ws.Cells.Item(2, 1).Value = 80
$ws.Cells.Item(2, 2).Value = "三商美邦人壽"
$ws.Cells.Item(2, 3).Value = "世紀理財變額萬能終身壽險"
$ws.Cells.Item(2, 4).Value = "蔡煌瑯"
$ws.Cells.Item(2, 5).Value = "保單號碼:649700013402自96年1月11日起迄今20年期(可隨時終止）"

$ws.Cells.Item(3, 1).Value = 81
$ws.Cells.Item(3, 2).Value = "三商美邦人壽"
$ws.Cells.Item(3, 3).Value = "世紀理財變額萬能終身壽險"
$ws.Cells.Item(3, 4).Value = "王琴賀"
$ws.Cells.Item(3, 5).Value = "保單號碼:649700007213自96年1月11日起迄今20年期(可隨時終止）"

$ws.Cells.Item(4, 1).Value = 82
$ws.Cells.Item(4, 2).Value = "新光人壽"
$ws.Cells.Item(4, 3).Value = "美利外幣終生還本型保險"
$ws.Cells.Item(4, 4).Value = "王琴賀"
$ws.Cells.Item(4, 5).Value = "保單號碼：1025099898自101年1月18日起迄今6年期"

# Drop the now-redundant trailing row (old row 5), shrinking the used range
# from A1:E5 to A1:E4.
$ws.Rows.Item(5).Delete()
